$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "289.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.79%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.60%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.950"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07183"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.86%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.825"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.48%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.688"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.34%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.748"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.91%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8972"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.28%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1654"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.91%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07708"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.55%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07954"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03041"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.77%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001494"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.02%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005847"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.99%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.459"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.07%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.33%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3319"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.70%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.56%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.032"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.47%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "20.08%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04509"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.34%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001216"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004617"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.10%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.17%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01565"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-7.96%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04357"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-6.83%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007305"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.99%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009886"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.58%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002073"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.01%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009514"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.87%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006049"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.14%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "173.87%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
